$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.07173050923138437
$ws.Range("J2").Value = 0.07173050923138438
$ws.Range("M2").Value = 64.88963566666666
$ws.Range("N2").Value = 194.668907
$ws.Range("O2").Value = 0.3123366816504561
$ws.Range("P2").Value = 0.3123366816504561
$ws.Range("Q2").Value = 1.994491101607778
$ws.Range("R2").Value = 17.95041991447
$ws.Range("S2").Value = 0.022404069226428
$ws.Range("T2").Value = 0.022404069226428
$ws.Range("I3").Value = 0.07173050923138437
$ws.Range("J3").Value = 0.07173050923138438
$ws.Range("O3").Value = 0.2703686564069002
$ws.Range("P3").Value = 0.2703686564069002
$ws.Range("S3").Value = 0.01939368140427214
$ws.Range("T3").Value = 0.01939368140427214
$ws.Range("I4").Value = 0.07173050923138437
$ws.Range("J4").Value = 0.07173050923138438
$ws.Range("M4").Value = 15.632391
$ws.Range("N4").Value = 46.897173
$ws.Range("O4").Value = 0.07524420627484885
$ws.Range("P4").Value = 0.07524420627484883
$ws.Range("Q4").Value = 0.48048759137
$ws.Range("R4").Value = 4.32438832233
$ws.Range("S4").Value = 0.005397305232806235
$ws.Range("T4").Value = 0.005397305232806235
$ws.Range("I5").Value = 0.07173050923138437
$ws.Range("J5").Value = 0.07173050923138438
$ws.Range("M5").Value = 71.062833
$ws.Range("N5").Value = 213.188499
$ws.Range("O5").Value = 0.3420504556677949
$ws.Range("P5").Value = 0.3420504556677948
$ws.Range("Q5").Value = 2.18423461031
$ws.Range("R5").Value = 19.65811149279
$ws.Range("S5").Value = 0.02453545336787799
$ws.Range("T5").Value = 0.02453545336787799
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3977653333333333
$ws.Range("H6").Value = 1.193296
$ws.Range("I6").Value = 0.9282694907686155
$ws.Range("J6").Value = 0.9282694907686156
$ws.Range("M6").Value = 64.88963566666666
$ws.Range("N6").Value = 194.668907
$ws.Range("O6").Value = 0.3123366816504561
$ws.Range("P6").Value = 0.3123366816504561
$ws.Range("Q6").Value = 25.81084756083022
$ws.Range("R6").Value = 232.297628047472
$ws.Range("S6").Value = 0.2899326124240281
$ws.Range("T6").Value = 0.2899326124240281
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3977653333333333
$ws.Range("H7").Value = 1.193296
$ws.Range("I7").Value = 0.9282694907686155
$ws.Range("J7").Value = 0.9282694907686156
$ws.Range("O7").Value = 0.2703686564069002
$ws.Range("P7").Value = 0.2703686564069002
$ws.Range("Q7").Value = 22.34269807462044
$ws.Range("R7").Value = 201.084282671584
$ws.Range("S7").Value = 0.2509749750026281
$ws.Range("T7").Value = 0.2509749750026281
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3977653333333333
$ws.Range("H8").Value = 1.193296
$ws.Range("I8").Value = 0.9282694907686155
$ws.Range("J8").Value = 0.9282694907686156
$ws.Range("M8").Value = 15.632391
$ws.Range("N8").Value = 46.897173
$ws.Range("O8").Value = 0.07524420627484885
$ws.Range("P8").Value = 0.07524420627484883
$ws.Range("Q8").Value = 6.218023216912
$ws.Range("R8").Value = 55.962208952208
$ws.Range("S8").Value = 0.06984690104204261
$ws.Range("T8").Value = 0.06984690104204259
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3977653333333333
$ws.Range("H9").Value = 1.193296
$ws.Range("I9").Value = 0.9282694907686155
$ws.Range("J9").Value = 0.9282694907686156
$ws.Range("M9").Value = 71.062833
$ws.Range("N9").Value = 213.188499
$ws.Range("O9").Value = 0.3420504556677949
$ws.Range("P9").Value = 0.3420504556677948
$ws.Range("Q9").Value = 28.266331455856
$ws.Range("R9").Value = 254.396983102704
$ws.Range("S9").Value = 0.3175150022999168
$ws.Range("T9").Value = 0.3175150022999168